$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = "LR1_signalling"
$ws.Cells.Item(2, 5).Value = 0.6077384310987872
$ws.Cells.Item(3, 5).Value = 0.5676121371873883
$ws.Cells.Item(4, 5).Value = 0.5807473489652029
$ws.Cells.Item(5, 3).Value = "LR3_other"
$ws.Cells.Item(5, 5).Value = 0.577983822890715
$ws.Cells.Item(6, 5).Value = 0.6875973831561824
$ws.Cells.Item(7, 5).Value = 0.6585521373376819
$ws.Cells.Item(8, 5).Value = 0.6745876697238239
$ws.Cells.Item(9, 5).Value = 0.7029086160885609
$ws.Cells.Item(10, 3).Value = "LR3_other"
$ws.Cells.Item(10, 5).Value = 0.5690872017843593
$ws.Cells.Item(11, 5).Value = 0.5604340364267851
$ws.Cells.Item(12, 3).Value = "LR1_signalling"
$ws.Cells.Item(12, 5).Value = 0.6246528812023089
$ws.Cells.Item(13, 3).Value = "LR3_other"
$ws.Cells.Item(13, 5).Value = 0.6566125041370295
$ws.Cells.Item(14, 5).Value = 0.5970461390090313
$ws.Cells.Item(15, 3).Value = "LR1_signalling"
$ws.Cells.Item(15, 5).Value = 0.6516624992957248
$ws.Cells.Item(16, 3).Value = "LR1_signalling"
$ws.Cells.Item(16, 5).Value = 0.5863605729765066
$ws.Cells.Item(17, 5).Value = 0.6592956398669667
$ws.Cells.Item(18, 5).Value = 0.6140413528525036
$ws.Cells.Item(19, 5).Value = 0.5739420476360719
$ws.Cells.Item(20, 5).Value = 0.5608710684452218
$ws.Cells.Item(21, 5).Value = 0.5436471462903746
$ws.Cells.Item(22, 5).Value = 0.5348007727792031
$ws.Cells.Item(23, 3).Value = "LR1_signalling"
$ws.Cells.Item(23, 5).Value = 0.7002801311176623
$ws.Cells.Item(24, 3).Value = "LR1_signalling"
$ws.Cells.Item(24, 5).Value = 0.6785881524847925
$ws.Cells.Item(25, 5).Value = 0.5647252059872219
$ws.Cells.Item(26, 5).Value = 0.6680321697468004
$ws.Cells.Item(27, 3).Value = "LR2_payload"
$ws.Cells.Item(27, 5).Value = 0.5968168985724113
$ws.Cells.Item(28, 5).Value = 0.6875826968663854
$ws.Cells.Item(29, 5).Value = 0.5718188593859584
$ws.Cells.Item(30, 5).Value = 0.5765003501019581
$ws.Cells.Item(31, 5).Value = 0.6655729117738267
$ws.Cells.Item(32, 5).Value = 0.6440185094176631
$ws.Cells.Item(33, 5).Value = 0.6110080358820285
$ws.Cells.Item(34, 3).Value = "LR3_other"
$ws.Cells.Item(34, 5).Value = 0.6815444062738008
$ws.Cells.Item(35, 3).Value = "LR1_signalling"
$ws.Cells.Item(35, 5).Value = 0.6834763353675209
$ws.Cells.Item(36, 5).Value = 0.5641521724472557
$ws.Cells.Item(37, 3).Value = "LR3_other"
$ws.Cells.Item(37, 5).Value = 0.666780003879725
$ws.Cells.Item(38, 5).Value = 0.5535057892274129
$ws.Cells.Item(39, 3).Value = "LR3_other"
$ws.Cells.Item(39, 5).Value = 0.6518979155981524
$ws.Cells.Item(40, 5).Value = 0.6356216615246701
$ws.Cells.Item(41, 5).Value = 0.6256987787171775
$ws.Cells.Item(42, 5).Value = 0.6920849156503396
$ws.Cells.Item(43, 5).Value = 0.6846679969513726
$ws.Cells.Item(44, 5).Value = 0.5552099764744204
$ws.Cells.Item(45, 5).Value = 0.6867315462096182
$ws.Cells.Item(46, 5).Value = 0.6226865745606156
$ws.Cells.Item(47, 5).Value = 0.6918521239613371
$ws.Cells.Item(48, 5).Value = 0.6240995449824275
$ws.Cells.Item(49, 3).Value = "LR1_signalling"
$ws.Cells.Item(49, 5).Value = 0.6986081537989165
$ws.Cells.Item(50, 3).Value = "LR2_payload"
$ws.Cells.Item(50, 5).Value = 0.5410678101919424
$ws.Cells.Item(51, 5).Value = 0.6949363903702647
$ws.Cells.Item(52, 5).Value = 0.6659942091969222
$ws.Cells.Item(53, 5).Value = 0.587471636619973
$ws.Cells.Item(54, 5).Value = 0.5691297066792967
$ws.Cells.Item(55, 5).Value = 0.5614261326668378
$ws.Cells.Item(56, 5).Value = 0.6877941059193086
$ws.Cells.Item(57, 5).Value = 0.7256233846173087
$ws.Cells.Item(58, 5).Value = 0.6974494397345261
$ws.Cells.Item(59, 5).Value = 0.6619611095095291
$ws.Cells.Item(60, 3).Value = "LR1_signalling"
$ws.Cells.Item(60, 5).Value = 0.6624516462878558
$ws.Cells.Item(61, 5).Value = 0.6785698847089061
$ws.Cells.Item(62, 5).Value = 0.5721511379575926
$ws.Cells.Item(63, 5).Value = 0.6669702385544822
$ws.Cells.Item(64, 5).Value = 0.567710888433699
$ws.Cells.Item(65, 5).Value = 0.5241378914469227
$ws.Cells.Item(66, 3).Value = "LR3_other"
$ws.Cells.Item(66, 5).Value = 0.6309033577831158
$ws.Cells.Item(67, 5).Value = 0.5270145400899887
$ws.Cells.Item(68, 3).Value = "LR2_payload"
$ws.Cells.Item(68, 5).Value = 0.5734978827351047
$ws.Cells.Item(69, 5).Value = 0.6621156476400103
$ws.Cells.Item(70, 5).Value = 0.5549423931360536
$ws.Cells.Item(71, 5).Value = 0.5592027956857066
$ws.Cells.Item(72, 5).Value = 0.6079451639677238
$ws.Cells.Item(73, 5).Value = 0.5937733520479671
$ws.Cells.Item(74, 5).Value = 0.6771617282013017
$ws.Cells.Item(75, 5).Value = 0.6424575790978876
$ws.Cells.Item(76, 3).Value = "LR2_payload"
$ws.Cells.Item(76, 5).Value = 0.5539725274224939
$ws.Cells.Item(77, 5).Value = 0.6701871539113655
$ws.Cells.Item(78, 5).Value = 0.6334426128454733
$ws.Cells.Item(79, 5).Value = 0.5597900734724914
$ws.Cells.Item(80, 3).Value = "LR2_payload"
$ws.Cells.Item(80, 5).Value = 0.6310394252531588
$ws.Cells.Item(81, 5).Value = 0.6161758231208476
$ws.Cells.Item(82, 5).Value = 0.6146393793606026
$ws.Cells.Item(83, 5).Value = 0.6445522025907785
$ws.Cells.Item(84, 5).Value = 0.6036984693615891
$ws.Cells.Item(85, 3).Value = "LR2_payload"
$ws.Cells.Item(85, 5).Value = 0.5746635055979473
$ws.Cells.Item(86, 5).Value = 0.5748408061268547
$ws.Cells.Item(87, 5).Value = 0.6566415377713357
$ws.Cells.Item(88, 5).Value = 0.6611246413887588
